$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("Jack", 1042.0),
    @("test", 2084.0),
    @("jg",   1042.0),
    @("jg",   3647.0),
    @("a",    521.0),
    @("jg",   521.0),
    @("jh",   521.0),
    @("jk",   4689.0),
    @("jk",   4689.0)
)

$row = 6
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
